# Applies the changes described by the commit:
# "added statistics about dead units and morale modifier"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 updates (existing units) ---
# Army 1 (columns B:K)
$ws.Range("F4").Value = 90          # Уязвимость урон: 10 -> 90
$ws.Range("I4").Value = "к+в"       # Кому наносит урон?: к -> к+в
$ws.Range("K4").Value = 15          # Количество отрядов: 3 -> 15

# Army 2 (columns N:W)
$ws.Range("T4").Value = -90         # Мораль: -50 -> -90
$ws.Range("W4").Value = 2           # Количество отрядов: 3 -> 2

# --- Legend update ---
$ws.Range("Z10").Value = "одиночный"   # особый -> одиночный

# --- Row 5 (new unit entry for Army 2) ---
$ws.Range("N5").Value = "Геныч"
$ws.Range("O5").Value = "п"
$ws.Range("P5").Value = 25
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = -90
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 1

# --- Selection / view state ---
$ws.Range("T11").Select()
$excel.ActiveWindow.ScrollColumn = 6
